# Change " and other subpages for proof.)" into three separate runs:
#   " and subpage" + "s" + ".)"
# (mirrors the author's edit which split the sentence into three runs)

$d = $word.ActiveDocument

$oldText = " and other subpages for proof.)"
$part1   = " and subpage"
$part2   = "s"
$part3   = ".)"

$rng = $d.Content
$found = $rng.Find.Execute($oldText)
if (-not $found) {
    throw "Could not find target text to replace"
}

$start = $rng.Start

# Remove the old text entirely.
$rng.Text = ""

# Re-insert the replacement text as three independent runs, matching the
# run-splitting seen in the target diff.
$pos = $start

$r1 = $d.Range($pos, $pos)
$r1.InsertAfter($part1)
$pos = $pos + $part1.Length

$r2 = $d.Range($pos, $pos)
$r2.InsertAfter($part2)
$pos = $pos + $part2.Length

$r3 = $d.Range($pos, $pos)
$r3.InsertAfter($part3)
$pos = $pos + $part3.Length
